$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 16.06960561832376
$ws.Cells.Item(2, 4).Value = 5.489965428507356
$ws.Cells.Item(2, 5).Value = 19.99195777004262
$ws.Cells.Item(2, 6).Value = 27.00008557441588
$ws.Cells.Item(2, 7).Value = 3.649689992866735
$ws.Cells.Item(2, 11).Value = 10.61424483852671
$ws.Cells.Item(2, 12).Value = 8.257987688500972
$ws.Cells.Item(2, 13).Value = 14.94732392085699
$ws.Cells.Item(2, 15).Value = 24.17652224851807
$ws.Cells.Item(3, 2).Value = 15.99587240615722
$ws.Cells.Item(3, 4).Value = 5.449532862029759
$ws.Cells.Item(3, 5).Value = 20.05554877090843
$ws.Cells.Item(3, 6).Value = 27.0161225467256
$ws.Cells.Item(3, 7).Value = 3.651601964372954
$ws.Cells.Item(3, 11).Value = 10.3388240298597
$ws.Cells.Item(3, 12).Value = 8.240403995425645
$ws.Cells.Item(3, 13).Value = 14.93224943701961
$ws.Cells.Item(3, 15).Value = 24.23659452891055
$ws.Cells.Item(4, 2).Value = 15.95350820909785
$ws.Cells.Item(4, 4).Value = 5.424161734941042
$ws.Cells.Item(4, 5).Value = 20.09674841339667
$ws.Cells.Item(4, 6).Value = 27.03296079563964
$ws.Cells.Item(4, 7).Value = 3.652838817959486
$ws.Cells.Item(4, 11).Value = 10.16419110774619
$ws.Cells.Item(4, 12).Value = 8.230738048883163
$ws.Cells.Item(4, 13).Value = 14.92512314534026
$ws.Cells.Item(4, 15).Value = 24.27869603626746
$ws.Cells.Item(5, 2).Value = 15.93698975145193
$ws.Cells.Item(5, 4).Value = 5.413689269901233
$ws.Cells.Item(5, 5).Value = 20.11408078905453
$ws.Cells.Item(5, 6).Value = 27.04157896035798
$ws.Cells.Item(5, 7).Value = 3.653358710883375
$ws.Cells.Item(5, 11).Value = 10.09170012450013
$ws.Cells.Item(5, 12).Value = 8.227086125390514
$ws.Cells.Item(5, 13).Value = 14.92275759979172
$ws.Cells.Item(5, 15).Value = 24.29716186997769
$ws.Cells.Item(6, 2).Value = 15.93429227835597
$ws.Cells.Item(6, 4).Value = 5.411942366641884
$ws.Cells.Item(6, 5).Value = 20.11699166574055
$ws.Cells.Item(6, 6).Value = 27.04311601856584
$ws.Cells.Item(6, 7).Value = 3.653445998373817
$ws.Cells.Item(6, 11).Value = 10.07958483656189
$ws.Cells.Item(6, 12).Value = 8.226497138606868
$ws.Cells.Item(6, 13).Value = 14.92239740772105
$ws.Cells.Item(6, 15).Value = 24.30030709254221
$ws.Cells.Item(7, 2).Value = 15.95328239995063
$ws.Cells.Item(7, 4).Value = 5.424021035466962
$ws.Cells.Item(7, 5).Value = 20.09697996264054
$ws.Cells.Item(7, 6).Value = 27.03306991433585
$ws.Cells.Item(7, 7).Value = 3.652845765114539
$ws.Cells.Item(7, 11).Value = 10.16321875486409
$ws.Cells.Item(7, 12).Value = 8.230687632110179
$ws.Cells.Item(7, 13).Value = 14.92508905873581
$ws.Cells.Item(7, 15).Value = 24.27893977609877
$ws.Cells.Item(8, 2).Value = 16.04358745808596
$ws.Cells.Item(8, 4).Value = 5.476139623666732
$ws.Cells.Item(8, 5).Value = 20.01343763308261
$ws.Cells.Item(8, 6).Value = 27.00416334442432
$ws.Cells.Item(8, 7).Value = 3.650336216553197
$ws.Cells.Item(8, 11).Value = 10.52046203494175
$ws.Cells.Item(8, 12).Value = 8.251691900580097
$ws.Cells.Item(8, 13).Value = 14.94168613251671
$ws.Cells.Item(8, 15).Value = 24.19615084519073
$ws.Cells.Item(9, 2).Value = 16.24309566797308
$ws.Cells.Item(9, 4).Value = 5.573861584000065
$ws.Cells.Item(9, 5).Value = 19.86664169152753
$ws.Cells.Item(9, 6).Value = 27.00298472004146
$ws.Cells.Item(9, 7).Value = 3.645911787766255
$ws.Cells.Item(9, 11).Value = 11.17469332811047
$ws.Cells.Item(9, 12).Value = 8.301722895108728
$ws.Cells.Item(9, 13).Value = 14.99098847548082
$ws.Cells.Item(9, 15).Value = 24.07530854811443
$ws.Cells.Item(10, 2).Value = 16.40240782205559
$ws.Cells.Item(10, 4).Value = 5.642722723291135
$ws.Cells.Item(10, 5).Value = 19.76908209347318
$ws.Cells.Item(10, 6).Value = 27.03594357769093
$ws.Cells.Item(10, 7).Value = 3.6429608554369
$ws.Cells.Item(10, 11).Value = 11.62410590861373
$ws.Cells.Item(10, 12).Value = 8.3436879352435
$ws.Cells.Item(10, 13).Value = 15.03722244657734
$ws.Cells.Item(10, 15).Value = 24.01197171746167
$ws.Cells.Item(11, 2).Value = 16.47743299916939
$ws.Cells.Item(11, 4).Value = 5.673370132035898
$ws.Cells.Item(11, 5).Value = 19.72691513951083
$ws.Cells.Item(11, 6).Value = 27.05825773411543
$ws.Cells.Item(11, 7).Value = 3.641682804069076
$ws.Cells.Item(11, 11).Value = 11.82122223161056
$ws.Cells.Item(11, 12).Value = 8.363866841860879
$ws.Cells.Item(11, 13).Value = 15.06038071077968
$ws.Cells.Item(11, 15).Value = 23.98870903653292
$ws.Cells.Item(12, 2).Value = 16.50619109385529
$ws.Cells.Item(12, 4).Value = 5.684874517637495
$ws.Cells.Item(12, 5).Value = 19.71126445556815
$ws.Cells.Item(12, 6).Value = 27.0677563572706
$ws.Cells.Item(12, 7).Value = 3.641208040752569
$ws.Cells.Item(12, 11).Value = 11.89477178655833
$ws.Cells.Item(12, 12).Value = 8.371660706813785
$ws.Cells.Item(12, 13).Value = 15.06945129442012
$ws.Cells.Item(12, 15).Value = 23.98069971286947
$ws.Cells.Item(13, 2).Value = 16.4999823759533
$ws.Cells.Item(13, 4).Value = 5.682401397148634
$ws.Cells.Item(13, 5).Value = 19.71462103042676
$ws.Cells.Item(13, 6).Value = 27.06566409158013
$ws.Cells.Item(13, 7).Value = 3.641309880762074
$ws.Cells.Item(13, 11).Value = 11.87898084494685
$ws.Cells.Item(13, 12).Value = 8.369975439566634
$ws.Cells.Item(13, 13).Value = 15.06748447008438
$ws.Cells.Item(13, 15).Value = 23.98238907056232
$ws.Cells.Item(14, 2).Value = 16.47979209022289
$ws.Cells.Item(14, 4).Value = 5.674318652057846
$ws.Cells.Item(14, 5).Value = 19.72562120179586
$ws.Cells.Item(14, 6).Value = 27.059018200215
$ws.Cells.Item(14, 7).Value = 3.641643560742425
$ws.Cells.Item(14, 11).Value = 11.82729534595459
$ws.Cells.Item(14, 12).Value = 8.364505014085131
$ws.Cells.Item(14, 13).Value = 15.06112094743925
$ws.Cells.Item(14, 15).Value = 23.98803406763415
$ws.Cells.Item(15, 2).Value = 16.46746964916573
$ws.Cells.Item(15, 4).Value = 5.66935446106848
$ws.Cells.Item(15, 5).Value = 19.73240037349471
$ws.Cells.Item(15, 6).Value = 27.05508383841571
$ws.Cells.Item(15, 7).Value = 3.641849147046808
$ws.Cells.Item(15, 11).Value = 11.79549285149755
$ws.Cells.Item(15, 12).Value = 8.361173965035533
$ws.Cells.Item(15, 13).Value = 15.05726216814543
$ws.Cells.Item(15, 15).Value = 23.99159598856398
$ws.Cells.Item(16, 2).Value = 16.3975544774483
$ws.Cells.Item(16, 4).Value = 5.640705889361009
$ws.Cells.Item(16, 5).Value = 19.77188223448398
$ws.Cells.Item(16, 6).Value = 27.03463229424246
$ws.Cells.Item(16, 7).Value = 3.643045669947981
$ws.Cells.Item(16, 11).Value = 11.61107292603508
$ws.Cells.Item(16, 12).Value = 8.342390766197218
$ws.Cells.Item(16, 13).Value = 15.03575138889448
$ws.Cells.Item(16, 15).Value = 24.01360381410601
$ws.Cells.Item(17, 2).Value = 16.35530407135228
$ws.Cells.Item(17, 4).Value = 5.622954684498573
$ws.Cells.Item(17, 5).Value = 19.79666908233567
$ws.Cells.Item(17, 6).Value = 27.02395851088659
$ws.Cells.Item(17, 7).Value = 3.643796145254993
$ws.Cells.Item(17, 11).Value = 11.49603104637447
$ws.Cells.Item(17, 12).Value = 8.331143946474151
$ws.Cells.Item(17, 13).Value = 15.02309653728487
$ws.Cells.Item(17, 15).Value = 24.02852759951986
$ws.Cells.Item(18, 2).Value = 16.33124428604005
$ws.Cells.Item(18, 4).Value = 5.612681070693036
$ws.Cells.Item(18, 5).Value = 19.81113423130384
$ws.Cells.Item(18, 6).Value = 27.01850863065424
$ws.Cells.Item(18, 7).Value = 3.644233857717029
$ws.Cells.Item(18, 11).Value = 11.42917484184096
$ws.Cells.Item(18, 12).Value = 8.324777807744987
$ws.Cells.Item(18, 13).Value = 15.01601822688923
$ws.Cells.Item(18, 15).Value = 24.03763363173513
$ws.Cells.Item(19, 2).Value = 16.323140123108
$ws.Cells.Item(19, 4).Value = 5.609191792103887
$ws.Cells.Item(19, 5).Value = 19.81606771529568
$ws.Cells.Item(19, 6).Value = 27.01678190397906
$ws.Cells.Item(19, 7).Value = 3.644383101648598
$ws.Cells.Item(19, 11).Value = 11.40642177486768
$ws.Cells.Item(19, 12).Value = 8.322640106882282
$ws.Cells.Item(19, 13).Value = 15.01365619459634
$ws.Cells.Item(19, 15).Value = 24.04080642001848
$ws.Cells.Item(20, 2).Value = 16.35977684091456
$ws.Cells.Item(20, 4).Value = 5.624850938534634
$ws.Cells.Item(20, 5).Value = 19.79400891960266
$ws.Cells.Item(20, 6).Value = 27.02502342886527
$ws.Cells.Item(20, 7).Value = 3.643715629131695
$ws.Cells.Item(20, 11).Value = 11.50834887879225
$ws.Cells.Item(20, 12).Value = 8.332330585553731
$ws.Cells.Item(20, 13).Value = 15.02442295771755
$ws.Cells.Item(20, 15).Value = 24.02688487109769
$ws.Cells.Item(21, 2).Value = 16.48571318830869
$ws.Cells.Item(21, 4).Value = 5.676695523383597
$ws.Cells.Item(21, 5).Value = 19.72238159076426
$ws.Cells.Item(21, 6).Value = 27.06094183463517
$ws.Cells.Item(21, 7).Value = 3.641545301328228
$ws.Cells.Item(21, 11).Value = 11.84250663058664
$ws.Cells.Item(21, 12).Value = 8.366107703320216
$ws.Cells.Item(21, 13).Value = 15.06298193653284
$ws.Cells.Item(21, 15).Value = 23.98635427787389
$ws.Cells.Item(22, 2).Value = 16.57003706149764
$ws.Cells.Item(22, 4).Value = 5.709987893065671
$ws.Cells.Item(22, 5).Value = 19.67741628939232
$ws.Cells.Item(22, 6).Value = 27.09052679199828
$ws.Cells.Item(22, 7).Value = 3.640180511281411
$ws.Cells.Item(22, 11).Value = 12.05450446003604
$ws.Cells.Item(22, 12).Value = 8.389070204987249
$ws.Cells.Item(22, 13).Value = 15.0899348582257
$ws.Cells.Item(22, 15).Value = 23.96452733418437
$ws.Cells.Item(23, 2).Value = 16.52485382152488
$ws.Cells.Item(23, 4).Value = 5.692274397688562
$ws.Cells.Item(23, 5).Value = 19.70124650088542
$ws.Cells.Item(23, 6).Value = 27.0741792690541
$ws.Cells.Item(23, 7).Value = 3.640904031825513
$ws.Cells.Item(23, 11).Value = 11.94195475061272
$ws.Cells.Item(23, 12).Value = 8.376734877069886
$ws.Cells.Item(23, 13).Value = 15.07539084925948
$ws.Cells.Item(23, 15).Value = 23.97574970317759
$ws.Cells.Item(24, 2).Value = 16.35775398222737
$ws.Cells.Item(24, 4).Value = 5.623993854260383
$ws.Cells.Item(24, 5).Value = 19.79521090981549
$ws.Cells.Item(24, 6).Value = 27.02453984036512
$ws.Cells.Item(24, 7).Value = 3.643752010986175
$ws.Cells.Item(24, 11).Value = 11.50278221723841
$ws.Cells.Item(24, 12).Value = 8.33179379483488
$ws.Cells.Item(24, 13).Value = 15.02382266864053
$ws.Cells.Item(24, 15).Value = 24.02762590980723
$ws.Cells.Item(25, 2).Value = 16.18681644787868
$ws.Cells.Item(25, 4).Value = 5.547926169633889
$ws.Cells.Item(25, 5).Value = 19.90454010248413
$ws.Cells.Item(25, 6).Value = 26.99735844433715
$ws.Cells.Item(25, 7).Value = 3.647055856390789
$ws.Cells.Item(25, 11).Value = 11.00299413224543
$ws.Cells.Item(25, 12).Value = 8.287259732484056
$ws.Cells.Item(25, 13).Value = 14.97587697182565
$ws.Cells.Item(25, 15).Value = 24.10354032586395
